$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4, 5, 6, 7, 8, 9, 10)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
